# B6-PowerPoint.pptx edit — Thu, May 14, 2020  2:04:53 AM
#
# Re-apply the custom table style ("Table_0",
# {E7BA431B-5640-42DD-B86C-4D77A6421EC2}) on the three tables in the deck
# with the built-in PowerPoint table style
# {7373AC3C-35EE-47A9-95ED-E37E8BCCA012}.
#
# The three tables each live as the first shape (a graphicFrame) on
# slides 14, 15 and 16.

$p = $ppt.ActivePresentation

$newStyleId = "{7373AC3C-35EE-47A9-95ED-E37E8BCCA012}"

$tableSlideIndexes = @(14, 15, 16)

foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
